$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_val data (regenerated to filter save games)
# Columns: B=TB, C=d2S, D=K, E=IP, G=sum (F=Win unchanged)
$rowNums = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15)
$colB    = @(0.1554434735375247,1.505614041169197,0.7287194209349384,3.182878228561681,3.182878228561681,1.505614041169197,1.505614041169197,0.02258322285507441,3.182878228561681,3.182878228561681,0.3464964993005633,0.7287194209349384,3.182878228561681,1.505614041169197)
$colC    = @(0.05231270169004087,0.3375848360084654,0.3375848360084654,1.65323645889881,1.65323645889881,1.65323645889881,1.65323645889881,1.65323645889881,1.65323645889881,1.65323645889881,1.65323645889881,0.3375848360084654,1.65323645889881,1.65323645889881)
$colD    = @(0.1529057820181812,0.7127328510149897,0.7127328510149897,0.1529057820181812,0.7127328510149897,3.082599426703578,3.082599426703578,0.7127328510149897,3.082599426703578,16.98373111632243,0.1529057820181812,3.082599426703578,0.1529057820181812,0.7127328510149897)
$colE    = @(0.4998867070740569,0.4998867070740569,0.4998867070740569,0.4998867070740569,6.48142807727062,0.4998867070740569,0.4998867070740569,0.4998867070740569,0.4998867070740569,0.4998867070740569,0.4998867070740569,0.4998867070740569,0.4998867070740569,0.4998867070740569)
$colG    = @(0.8605486643198037,3.055818435266709,2.27892381503245,5.488907176552729,12.0302756157461,6.741336633845642,6.741336633845642,2.888439239842931,8.418600821238126,22.31973251085698,2.652525447291612,4.64879039072104,5.488907176552729,4.371470058157054)

for ($i = 0; $i -lt $rowNums.Length; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]   # B: TB
    $ws.Cells.Item($r, 3).Value = $colC[$i]   # C: d2S
    $ws.Cells.Item($r, 4).Value = $colD[$i]   # D: K
    $ws.Cells.Item($r, 5).Value = $colE[$i]   # E: IP
    $ws.Cells.Item($r, 7).Value = $colG[$i]   # G: sum
}
